# Implemented hamilton breath detection
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the filename referenced by the shared string table (A2)
$ws.Range("A2").Value = "Waves_009.txt"

# Update the detection/measurement results in row 2 (columns D through AK)
$ws.Range("D2").Value = 72
$ws.Range("E2").Value = 17
$ws.Range("F2").Value = 452.55
$ws.Range("G2").Value = 7.54
$ws.Range("H2").Value = 4.67
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4.19
$ws.Range("L2").Value = 3.13
$ws.Range("M2").Value = 0.09
$ws.Range("N2").Value = 2.13
$ws.Range("O2").Value = 1.16
$ws.Range("P2").Value = 0.03
$ws.Range("Q2").Value = 56.8
$ws.Range("R2").Value = 33.51
$ws.Range("S2").Value = 0.93
$ws.Range("T2").Value = 2.12
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 94.87
$ws.Range("X2").Value = 38.18
$ws.Range("Y2").Value = 1.05
$ws.Range("Z2").Value = 6.56
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 3.89
$ws.Range("AD2").Value = 2.53
$ws.Range("AE2").Value = 0.07000000000000001
$ws.Range("AF2").Value = 7.83
$ws.Range("AG2").Value = 2.99
$ws.Range("AH2").Value = 0.08
$ws.Range("AI2").Value = 11.46
$ws.Range("AJ2").Value = 3.4
$ws.Range("AK2").Value = 0.09
